$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 9: BIM 5D - Custos e Orçamentos -> mark AP/RP column as "AP" ---
$ws.Range("L9").Value = "AP"

# --- Row 10: BIM 6D - Análise Energética -> record the "prova do módulo energético" ---
$ws.Range("J10").Value = 48
$ws.Range("L10").Value = "AP"

# M10/M11 pick up the same highlighted "completed" look already used by M4:M9,
# so copy that formatting down before filling in the next sequence numbers.
$ws.Range("M9").Copy()
$ws.Range("M10").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("M9").Copy()
$ws.Range("M11").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("M10").Value = 8

# --- Row 11: BIM 7D - Gestão de Instalações -> next sequence marker ---
$ws.Range("M11").Value = 9

# --- Column widths: narrow column C, add a small spacer column E ---
$ws.Columns.Item(3).ColumnWidth = 19
$ws.Columns.Item(5).ColumnWidth = 4.85546875

# --- Selection moves to E17 ---
$ws.Range("E17").Select()
